$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.436.30'
$ws.Range("E2").Value = '  +1.90%  '
$ws.Range("D3").Value = '3.785.31'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '626.11'
$ws.Range("E5").Value = '  +4.08%  '
$ws.Range("D6").Value = '164.24'
$ws.Range("E6").Value = '  -0.52%  '
$ws.Range("D7").Value = '3.780.48'
$ws.Range("E7").Value = '  -0.38%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").Value = '0.522'
$ws.Range("E9").Value = '  +0.85%  '
$ws.Range("E10").Value = '  +1.37%  '
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("D12").Value = '6.62'
$ws.Range("E12").Value = '  +2.01%  '
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").Value = '35.50'
$ws.Range("E14").Value = '  -0.63%  '
$ws.Range("D15").Value = '4.424.08'
$ws.Range("E15").Value = '  -0.31%  '
$ws.Range("D16").Value = '3.667.15'
$ws.Range("E16").Value = '  -3.76%  '
$ws.Range("D17").Value = '69.467.84'
$ws.Range("E17").Value = '  +1.90%  '
$ws.Range("D18").Value = '17.89'
$ws.Range("E18").Value = '  -2.34%  '
$ws.Range("D19").Value = '7.11'
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("E20").Value = '  -1.32%  '
$ws.Range("D21").Value = '469.27'
$ws.Range("E21").Value = '  +1.57%  '
$ws.Range("D22").Value = '9.63'
$ws.Range("E22").Value = '  -0.72%  '
$ws.Range("D23").Value = '0.703'
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").Value = '0.0000150'
$ws.Range("E24").Value = '  +1.92%  '
$ws.Range("D25").Value = '83.26'
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").Value = '12.05'
$ws.Range("E26").Value = '  +0.37%  '
$ws.Range("B27").Value = 'Fetch.AI'
$ws.Range("C27").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D27").Value = '2.16'
$ws.Range("E27").Value = '  +2.27%  '
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '9.99'
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").Value = '3.935.28'
$ws.Range("E30").Value = '  -0.35%  '
$ws.Range("D31").Value = '2.67'
$ws.Range("E31").Value = '  +1.40%  '
$ws.Range("D32").Value = '2.24'
$ws.Range("E32").Value = '  -0.15%  '
$ws.Range("D33").Value = '7.31'
$ws.Range("E33").Value = '  -0.44%  '
$ws.Range("D34").Value = '28.89'
$ws.Range("E34").Value = '  -1.54%  '
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").Value = '9.02'
$ws.Range("E36").Value = '  -0.11%  '
$ws.Range("D37").Value = '3.738.04'
$ws.Range("E37").Value = '  -0.30%  '
$ws.Range("E38").Value = '  +2.60%  '
$ws.Range("E39").Value = '  +8.48%  '
$ws.Range("D40").Value = '3.34'
$ws.Range("E40").Value = '  +0.56%  '
$ws.Range("D41").Value = '5.81'
$ws.Range("E41").Value = '  -0.70%  '
$ws.Range("D42").Value = '0.968'
$ws.Range("E42").Value = '  -2.01%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("E46").Value = '  +3.44%  '
$ws.Range("D47").Value = '153.51'
$ws.Range("E47").Value = '  +0.82%  '
$ws.Range("D48").Value = '42.93'
$ws.Range("E48").Value = '  -0.38%  '
$ws.Range("D49").Value = '46.82'
$ws.Range("E49").Value = '  -1.59%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '8.42'
$ws.Range("E50").Value = '  +0.69%  '
$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D51").Value = '1.38'
$ws.Range("E51").Value = '  +1.84%  '
